$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 28.30556419040565
$ws.Range("A3").Value = 7.676473805426895
$ws.Range("A4").Value = 7.944418227584919
$ws.Range("A5").Value = 3.415618938336195
$ws.Range("A6").Value = 3.116020622501082
$ws.Range("A7").Value = 2.994456602645386
$ws.Range("A8").Value = 5.88973572244214
$ws.Range("A9").Value = 1.534354667387845
$ws.Range("A10").Value = 10.12321440681649
$ws.Range("A11").Value = 0.5900651275773328
$ws.Range("A12").Value = 6.028896973072136
$ws.Range("A13").Value = 3.609487352116929
$ws.Range("A14").Value = 4.609890264698635
$ws.Range("A15").Value = 0.8522022127903881
$ws.Range("A16").Value = 1.987579744537641
$ws.Range("A17").Value = 1.257736972979956
$ws.Range("A18").Value = 5.216982649957799
$ws.Range("A19").Value = 3.788602647592455
$ws.Range("A20").Value = 0.6387614721413115
$ws.Range("A21").Value = 3.414959208046184
$ws.Range("A22").Value = 2.038917012410053
$ws.Range("A23").Value = 0.4854494366466895
$ws.Range("A24").Value = 2.808907884265096
$ws.Range("A25").Value = 1.293759942376454
$ws.Range("A26").Value = 5.710266940505818
$ws.Range("A27").Value = 5.322634715292821
$ws.Range("A28").Value = 1.271915641835164
$ws.Range("A29").Value = 4.956427830928106
$ws.Range("A30").Value = 3.458700048572354
$ws.Range("A31").Value = 6.300969794730321
$ws.Range("A32").Value = 3.947004899624091
$ws.Range("A33").Value = 4.448180521353947
$ws.Range("A34").Value = 3.040576669824759
$ws.Range("A35").Value = 1.1691830588328
$ws.Range("A36").Value = 4.192214256073981
$ws.Range("A37").Value = 3.828851996518921
$ws.Range("A38").Value = 6.237926996838752
$ws.Range("A39").Value = 5.348889142888623
$ws.Range("A40").Value = 6.490047176699989
$ws.Range("A41").Value = 6.966328863329551
$ws.Range("A42").Value = 8.692575250883891
$ws.Range("A43").Value = 1.65804762674108
$ws.Range("A44").Value = 4.869295015244063
$ws.Range("A45").Value = 0.1385215065000978
$ws.Range("A46").Value = 2.560196570020395
$ws.Range("A47").Value = 3.960542001562146
$ws.Range("A48").Value = 1.587850634691875
$ws.Range("A49").Value = 4.497191368304243
$ws.Range("A50").Value = 2.687779156573725
$ws.Range("A51").Value = 3.35918414536161
$ws.Range("A52").Value = 2.227878816095682
$ws.Range("A53").Value = 0.1157204711687996
$ws.Range("A54").Value = 2.86152046210637
$ws.Range("A55").Value = 2.700725619718668
$ws.Range("A56").Value = 3.261091502558912
$ws.Range("A57").Value = 0.04716473961630641
$ws.Range("A58").Value = 2.259787257745074
$ws.Range("A59").Value = 1.051048077765387
$ws.Range("A60").Value = 0.04896573817785566
$ws.Range("A61").Value = 1.246089703777358
$ws.Range("A62").Value = 3.732817596500496
$ws.Range("A63").Value = 7.065275994188823
$ws.Range("A64").Value = 6.788761644023396
$ws.Range("A65").Value = 0.9023699298346912
$ws.Range("A66").Value = 7.210383439699797
$ws.Range("A67").Value = 2.931111611324496
$ws.Range("A68").Value = 5.354712804239426
$ws.Range("A69").Value = 7.50997915616145
$ws.Range("A70").Value = 3.300823105730927
$ws.Range("A71").Value = 0.6462444035926467
$ws.Range("A72").Value = 4.247627865884226
$ws.Range("A73").Value = 2.414099424793022
$ws.Range("A74").Value = 6.812554763277262
$ws.Range("A75").Value = 0.8685221301803381
$ws.Range("A76").Value = 1.24962359301594
$ws.Range("A77").Value = 1.351711336541911
$ws.Range("A78").Value = 2.646163582967546
$ws.Range("A79").Value = 2.122528569092339
$ws.Range("A80").Value = 0.2238756736561527
$ws.Range("A81").Value = 3.240177278725497
$ws.Range("A82").Value = 0.3967112054478719
$ws.Range("A83").Value = 1.883477187976155
$ws.Range("A84").Value = 3.332656956713947
$ws.Range("A85").Value = 0.4007973560852918
$ws.Range("A86").Value = 2.987432433173382
